$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text looks like a plain number need an explicit Text
# number format first, otherwise Excel auto-converts the string into a
# numeric value (changing the stored cell type).

$ws.Range('D2').Value = '64.203.34'
$ws.Range('E2').Value = '  -2.33%  '
$ws.Range('D3').Value = '3.480.94'
$ws.Range('E3').Value = '  -3.32%  '
$ws.Range('E4').Value = '  +0.08%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '583.47'
$ws.Range('E5').Value = '  -3.45%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '131.50'
$ws.Range('E6').Value = '  -3.83%  '
$ws.Range('D7').Value = '3.480.63'
$ws.Range('E7').Value = '  -3.34%  '
$ws.Range('E8').Value = '  +0.03%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.489'
$ws.Range('E9').Value = '  -1.72%  '
$ws.Range('E10').Value = '  -1.41%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '7.23'
$ws.Range('E11').Value = '  -0.44%  '
$ws.Range('E12').Value = '  -1.50%  '
$ws.Range('D13').Value = '4.075.85'
$ws.Range('E13').Value = '  -3.22%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '27.75'
$ws.Range('E14').Value = '  -1.00%  '
$ws.Range('E15').Value = '  -4.54%  '
$ws.Range('E16').Value = '  +0.25%  '
$ws.Range('D17').Value = '3.482.29'
$ws.Range('E17').Value = '  -3.28%  '
$ws.Range('D18').Value = '64.266.58'
$ws.Range('E18').Value = '  -2.45%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '9.90'
$ws.Range('E19').Value = '  -1.57%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '14.26'
$ws.Range('E20').Value = '  -2.85%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '5.65'
$ws.Range('E21').Value = '  -4.11%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '393.59'
$ws.Range('E22').Value = '  -0.68%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.575'
$ws.Range('E23').Value = '  -2.63%  '
$ws.Range('D24').Value = '3.623.08'
$ws.Range('E24').Value = '  -3.31%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '73.12'
$ws.Range('E25').Value = '  -1.72%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.999'
$ws.Range('E26').Value = '  -0.14%  '
$ws.Range('E27').Value = '  -8.91%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '1.56'
$ws.Range('E28').Value = '  -6.27%  '
$ws.Range('B29').Value = 'RenderToken'
$ws.Range('C29').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '7.47'
$ws.Range('E29').Value = '  -9.07%  '
$ws.Range('B30').Value = 'Binance-PegBSC-USD'
$ws.Range('C30').Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.995'
$ws.Range('E30').Value = '  -0.44%  '
$ws.Range('E31').Value = '  -6.90%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '8.18'
$ws.Range('E32').Value = '  -5.20%  '
$ws.Range('D33').Value = '3.482.99'
$ws.Range('E33').Value = '  -3.30%  '
$ws.Range('E34').Value = '  +0.02%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '23.86'
$ws.Range('E35').Value = '  -2.73%  '
$ws.Range('E36').Value = '  -2.51%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '5.26'
$ws.Range('E37').Value = '  -2.81%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '6.98'
$ws.Range('E38').Value = '  -1.81%  '
$ws.Range('B39').Value = 'ImmutableX'
$ws.Range('C39').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.57'
$ws.Range('E39').Value = '  -2.50%  '
$ws.Range('B40').Value = 'Monero'
$ws.Range('C40').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '169.91'
$ws.Range('E40').Value = '  -0.11%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.0804'
$ws.Range('E41').Value = '  -4.01%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.810'
$ws.Range('E42').Value = '  -3.83%  '
$ws.Range('B43').Value = 'FirstDigitalUSD'
$ws.Range('C43').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.00'
$ws.Range('E43').Value = '  +0.04%  '
$ws.Range('B44').Value = 'EnergySwap'
$ws.Range('C44').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '25.42'
$ws.Range('E44').Value = '  -4.62%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '41.63'
$ws.Range('E45').Value = '  -3.76%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '1.19'
$ws.Range('E46').Value = '  -5.98%  '
$ws.Range('E47').Value = '  -4.40%  '
$ws.Range('E48').Value = '  -3.98%  '
$ws.Range('B49').Value = 'Cosmos'
$ws.Range('C49').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '6.87'
$ws.Range('E49').Value = '  -2.76%  '
$ws.Range('B50').Value = 'Maker'
$ws.Range('C50').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D50').Value = '2.426.66'
$ws.Range('E50').Value = '  -0.26%  '
$ws.Range('E51').Value = '  -1.99%  '
